# Update investment-cost result values for sheets 2025, 2030, 2035
# as supplied by the latest server run.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 3481.11133040007
$ws.Range("E2").Value = 288784.0257356465
$ws.Range("I2").Value = 157511.7087451186
$ws.Range("L2").Value = 489912.4167596999
$ws.Range("M2").Value = 112998.9730188
$ws.Range("N2").Value = 72668.9505637426
$ws.Range("O2").Value = 69732.90075344281

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 2209.181339077783
$ws.Range("B2").Value = 57019.8100017786
$ws.Range("E2").Value = 261026.8276704075
$ws.Range("I2").Value = 283850.8070745453
$ws.Range("L2").Value = 100844.7993015671
$ws.Range("M2").Value = 106564.8742708
$ws.Range("N2").Value = 33073.85645417624
$ws.Range("O2").Value = 22743.17640421139

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 23566.89669566984
$ws.Range("B2").Value = 12464.79578491448
$ws.Range("E2").Value = 109707.1026912166
$ws.Range("I2").Value = 158247.5988221456
$ws.Range("M2").Value = 62871.57057687733
$ws.Range("N2").Value = 48867.45187110166
$ws.Range("O2").Value = 59761.20443857631
